$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.381382609085108
$ws.Range("C2").Value = 0.05879030905613547
$ws.Range("D2").Value = 0.3478214565051303
$ws.Range("F2").Value = 1.179065208724786
$ws.Range("G2").Value = 0.5527075265866586
$ws.Range("H2").Value = 0.6979076677046265
$ws.Range("J2").Value = 0.3486957219232778
$ws.Range("K2").Value = 0.363631907206269
$ws.Range("O2").Value = 2.474723568599416

$ws.Range("B3").Value = 0.3389749717559596
$ws.Range("C3").Value = 0.05353914119946523
$ws.Range("D3").Value = 0.3351820282697417
$ws.Range("F3").Value = 1.182757566074585
$ws.Range("G3").Value = 0.5577120479550857
$ws.Range("H3").Value = 0.7041546484010084
$ws.Range("J3").Value = 0.3372058456523774
$ws.Range("K3").Value = 0.3186676425652593
$ws.Range("O3").Value = 2.49815563862731

$ws.Range("B4").Value = 0.3129154250683825
$ws.Range("C4").Value = 0.0503003378455702
$ws.Range("D4").Value = 0.3275651915319742
$ws.Range("F4").Value = 1.185723802380089
$ws.Range("G4").Value = 0.5611861172651942
$ws.Range("H4").Value = 0.7083057156723243
$ws.Range("J4").Value = 0.3303527125874268
$ws.Range("K4").Value = 0.2909835266092671
$ws.Range("O4").Value = 2.51404521357982

$ws.Range("B5").Value = 0.3022913321089504
$ws.Range("C5").Value = 0.04897693632322841
$ws.Range("D5").Value = 0.3244976132567956
$ws.Range("F5").Value = 1.187108364809603
$ws.Range("G5").Value = 0.5627026188610671
$ws.Range("H5").Value = 0.7100766470375675
$ws.Range("J5").Value = 0.3276107333784353
$ws.Range("K5").Value = 0.2796836746927625
$ws.Range("O5").Value = 2.520897896684403

$ws.Range("B6").Value = 0.3005269510116761
$ws.Range("C6").Value = 0.04875697426346903
$ws.Range("D6").Value = 0.3239904446155748
$ws.Range("F6").Value = 1.187348888038038
$ws.Range("G6").Value = 0.562960516948678
$ws.Range("H6").Value = 0.7103755014995627
$ws.Range("J6").Value = 0.3271584939745367
$ws.Range("K6").Value = 0.2778062561407069
$ws.Range("O6").Value = 2.522058577289627

$ws.Range("B7").Value = 0.312772162615147
$ws.Range("C7").Value = 0.05028250428418346
$ws.Range("D7").Value = 0.3275236737024159
$ws.Range("F7").Value = 1.1857417632799
$ws.Range("G7").Value = 0.5612061613911337
$ws.Range("H7").Value = 0.7083292777990664
$ws.Range("J7").Value = 0.3303155278862562
$ws.Range("K7").Value = 0.2908312061828155
$ws.Range("O7").Value = 2.514136102814362

$ws.Range("B8").Value = 0.3667652300102588
$ws.Range("C8").Value = 0.05698278739802731
$ws.Range("D8").Value = 0.3434336339579431
$ws.Range("F8").Value = 1.180193229632465
$ws.Range("G8").Value = 0.5543497330687472
$ws.Range("H8").Value = 0.6999961837329138
$ws.Range("J8").Value = 0.3446921518173269
$ws.Range("K8").Value = 0.3481444000663032
$ws.Range("O8").Value = 2.482491055285024

$ws.Range("B9").Value = 0.4724532020235586
$ws.Range("C9").Value = 0.07000270127902297
$ws.Range("D9").Value = 0.3757681964598873
$ws.Range("F9").Value = 1.174860992600429
$ws.Range("G9").Value = 0.5440935633866033
$ws.Range("H9").Value = 0.6861562691596532
$ws.Range("J9").Value = 0.3744865127727905
$ws.Range("K9").Value = 0.459906691479631
$ws.Range("O9").Value = 2.432363357523712

$ws.Range("B10").Value = 0.5499590261626679
$ws.Range("C10").Value = 0.07949156352151476
$ws.Range("D10").Value = 0.4002111384830869
$ws.Range("F10").Value = 1.174329389647895
$ws.Range("G10").Value = 0.5385103590625491
$ws.Range("H10").Value = 0.6775112198893751
$ws.Range("J10").Value = 0.3973577333293719
$ws.Range("K10").Value = 0.5416072628156599
$ws.Range("O10").Value = 2.402819489974519

$ws.Range("B11").Value = 0.5851823495938504
$ws.Range("C11").Value = 0.08379078986015998
$ws.Range("D11").Value = 0.411479026979606
$ws.Range("F11").Value = 1.174823657823964
$ws.Range("G11").Value = 0.5363959576298356
$ws.Range("H11").Value = 0.6739087525737233
$ws.Range("J11").Value = 0.4079767040252023
$ws.Range("K11").Value = 0.5786802844975512
$ws.Range("O11").Value = 2.390964005041496

$ws.Range("B12").Value = 0.598514960558731
$ws.Range("C12").Value = 0.08541621911284381
$ws.Range("D12").Value = 0.4157671187782057
$ws.Range("F12").Value = 1.175116713587414
$ws.Range("G12").Value = 0.535656596757768
$ws.Range("H12").Value = 0.6725920530879677
$ws.Range("J12").Value = 0.4120287440035639
$ws.Range("K12").Value = 0.5927048750546362
$ws.Range("O12").Value = 2.386702682016335

$ws.Range("B13").Value = 0.5956438093934935
$ws.Range("C13").Value = 0.08506627091112762
$ws.Range("D13").Value = 0.4148426622474233
$ws.Range("F13").Value = 1.175048888845438
$ws.Range("G13").Value = 0.5358131020176415
$ws.Range("H13").Value = 0.6728735164499042
$ws.Range("J13").Value = 0.4111546921247964
$ws.Range("K13").Value = 0.5896850707463557
$ws.Range("O13").Value = 2.387610285124339

$ws.Range("B14").Value = 0.5862793499384509
$ws.Range("C14").Value = 0.08392456750618749
$ws.Range("D14").Value = 0.4118313873534305
$ws.Range("F14").Value = 1.174845645352775
$ws.Range("G14").Value = 0.5363339002745846
$ws.Range("H14").Value = 0.6737994754499326
$ws.Range("J14").Value = 0.4083094492646069
$ws.Range("K14").Value = 0.5798343842560598
$ws.Range("O14").Value = 2.390608850474592

$ws.Range("B15").Value = 0.5805425842355874
$ws.Range("C15").Value = 0.08322490019263284
$ws.Range("D15").Value = 0.4099896480025507
$ws.Range("F15").Value = 1.174734943561774
$ws.Range("G15").Value = 0.5366608937140143
$ws.Range("H15").Value = 0.674372834739124
$ws.Range("J15").Value = 0.4065706751553222
$ws.Range("K15").Value = 0.573798687206363
$ws.Range("O15").Value = 2.392475271210103

$ws.Range("B16").Value = 0.5476563422012646
$ws.Range("C16").Value = 0.07921024167347923
$ws.Range("D16").Value = 0.3994777307320021
$ws.Range("F16").Value = 1.174311903427096
$ws.Range("G16").Value = 0.5386571129313609
$ws.Range("H16").Value = 0.6777532936506958
$ws.Range("J16").Value = 0.3966680786244865
$ws.Range("K16").Value = 0.5391825152649972
$ws.Range("O16").Value = 2.40362617617123

$ws.Range("B17").Value = 0.5274723424518584
$ws.Range("C17").Value = 0.07674287088138954
$ws.Range("D17").Value = 0.3930669519477306
$ws.Range("F17").Value = 1.174240941004001
$ws.Range("G17").Value = 0.5399907914509114
$ws.Range("H17").Value = 0.6799116651333463
$ws.Range("J17").Value = 0.3906481471247929
$ws.Range("K17").Value = 0.517922239324804
$ws.Range("O17").Value = 2.410872845838966

$ws.Range("B18").Value = 0.5158598397972582
$ws.Range("C18").Value = 0.07532208196114709
$ws.Range("D18").Value = 0.3893936467115111
$ws.Range("F18").Value = 1.174269408072071
$ws.Range("G18").Value = 0.5407979228002802
$ws.Range("H18").Value = 0.6811841845528264
$ws.Range("J18").Value = 0.3872058547183315
$ws.Range("K18").Value = 0.5056851862775318
$ws.Range("O18").Value = 2.4151900307467

$ws.Range("B19").Value = 0.5119275189797179
$ws.Range("C19").Value = 0.07484075182087224
$ws.Range("D19").Value = 0.3881523399123523
$ws.Range("F19").Value = 1.174290944763442
$ws.Range("G19").Value = 0.5410780752928801
$ws.Range("H19").Value = 0.6816203762107023
$ws.Range("J19").Value = 0.3860438256243128
$ws.Range("K19").Value = 0.5015404594148265
$ws.Range("O19").Value = 2.416677353236764

$ws.Range("B20").Value = 0.5296213002316108
$ws.Range("C20").Value = 0.0770056954078342
$ws.Range("D20").Value = 0.3937479422006049
$ws.Range("F20").Value = 1.174241324292908
$ws.Range("G20").Value = 0.5398446746232395
$ws.Range("H20").Value = 0.6796786860636672
$ws.Range("J20").Value = 0.3912868872340312
$ws.Range("K20").Value = 0.5201863386091077
$ws.Range("O20").Value = 2.410085992889222

$ws.Range("B21").Value = 0.5890300791917582
$ws.Range("C21").Value = 0.0842599847260459
$ws.Range("D21").Value = 0.4127152983071198
$ws.Range("F21").Value = 1.174902468855649
$ws.Range("G21").Value = 0.5361792639536418
$ws.Range("H21").Value = 0.6735262103963322
$ws.Range("J21").Value = 0.4091443284689262
$ws.Range("K21").Value = 0.5827281601744687
$ws.Range("O21").Value = 2.389721906492696

$ws.Range("B22").Value = 0.6278236009427189
$ws.Range("C22").Value = 0.08898591091815433
$ws.Range("D22").Value = 0.4252349666077464
$ws.Range("F22").Value = 1.175951772769579
$ws.Range("G22").Value = 0.5341411738586288
$ws.Range("H22").Value = 0.6697819393632187
$ws.Range("J22").Value = 0.4209951039843105
$ws.Range("K22").Value = 0.6235199680217818
$ws.Range("O22").Value = 2.377742426100525

$ws.Range("B23").Value = 0.6071220823030217
$ws.Range("C23").Value = 0.0864650182089548
$ws.Range("D23").Value = 0.4185417559513382
$ws.Range("F23").Value = 1.175335252958547
$ws.Range("G23").Value = 0.535196186142656
$ws.Range("H23").Value = 0.6717550089410054
$ws.Range("J23").Value = 0.4146536652061883
$ws.Range("K23").Value = 0.6017564562402242
$ws.Range("O23").Value = 2.384014343754984

$ws.Range("B24").Value = 0.5286497826913035
$ws.Range("C24").Value = 0.07688687948576955
$ws.Range("D24").Value = 0.3934400280658963
$ws.Range("F24").Value = 1.174240935256691
$ws.Range("G24").Value = 0.5399106082646838
$ws.Range("H24").Value = 0.6797839173335518
$ws.Range("J24").Value = 0.3909980547511793
$ws.Range("K24").Value = 0.5191627836107671
$ws.Range("O24").Value = 2.410441258840294

$ws.Range("B25").Value = 0.4438850789687194
$ws.Range("C25").Value = 0.06649367299102948
$ws.Range("D25").Value = 0.3668998584732321
$ws.Range("F25").Value = 1.175709139576362
$ws.Range("G25").Value = 0.5465258931144419
$ws.Range("H25").Value = 0.6896327324470732
$ws.Range("J25").Value = 0.3662543629464636
$ws.Range("K25").Value = 0.4297422363461862
$ws.Range("O25").Value = 2.444645865868821
